# The template used to contain a hard-coded "M2Doc version mismatch" error
# banner (four spaces, an orange "<---" marker, the mismatch message, and
# four trailing spaces) inserted right after the spell-check boundary for
# "demonstration". Now that the test template version has been fixed, that
# banner text is obsolete and must be removed so the paragraph reads
# "A simple demonstration of a query :" again.

$d = $word.ActiveDocument

$bannerText = "    <---M2Doc version mismatch: template is 3.1.1 and runtime is 3.2.0    "

$searchRange = $d.Content
$searchRange.Find.ClearFormatting()

$found = $searchRange.Find.Execute($bannerText, $true, $false, $false, $false,
                                    $false, $true, 1, $false, "", 0)

if ($found -and $searchRange.Find.Found) {
    # Delete the whole matched range (all of the runs that make up the
    # banner), leaving the surrounding text/runs untouched.
    $searchRange.Delete()
} else {
    Write-Host "Version mismatch banner not found; no changes made."
}
